$d = $word.ActiveDocument

$replacements = @(
    @("411×2=", "541×3="),
    @("577×2=", "687×2="),
    @("736×3=", "704×5="),
    @("345×2=", "640×9="),
    @("391×9=", "220×2="),
    @("834×5=", "220×8="),
    @("582×9=", "707×4="),
    @("114×7=", "583×9="),
    @("669×3=", "780×4="),
    @("446×6=", "339×4="),
    @("711×5=", "461×5="),
    @("669×6=", "410×5="),
    @("571×7=", "253×7="),
    @("480×7=", "578×5="),
    @("278×2=", "168×4="),
    @("943×6=", "688×7="),
    @("777×8=", "302×8="),
    @("478×4=", "840×8="),
    @("556×8=", "208×6="),
    @("583×3=", "687×5="),
    @("210×6=", "766×9="),
    @("492×3=", "675×2="),
    @("446×5=", "332×8="),
    @("134×9=", "269×6="),
    @("663×5=", "328×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
